$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete row 108 ("「イルカのおおきなジャンプ！」") entirely, shifting all
# subsequent rows up by one.
$ws.Rows.Item(108).Delete()
